# Insert two new rows of Acelga price data at rows 193-194, pushing the
# existing rows 193:289 down to 195:291 (weekly update: newest entries on top).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 193, shifting rows 193:289 down.
$ws.Rows("193:194").Insert()

# New row 193: same as (what is now) row 195 except Fecha (D) and Volumen (J).
$ws.Range("A193").Value = $ws.Range("A195").Value()
$ws.Range("B193").Value = $ws.Range("B195").Value()
$ws.Range("C193").Value = $ws.Range("C195").Value()
$ws.Range("D193").Value = 44837
$ws.Range("D193").NumberFormat = $ws.Range("D195").NumberFormat
$ws.Range("E193").Value = $ws.Range("E195").Value()
$ws.Range("F193").Value = $ws.Range("F195").Value()
$ws.Range("G193").Value = $ws.Range("G195").Value()
$ws.Range("H193").Value = $ws.Range("H195").Value()
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 200
$ws.Range("K193").Value = 700
$ws.Range("L193").Value = 800
$ws.Range("M193").Value = 750
$ws.Range("N193").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O193").Value = "Provincia de Diguillín"
$ws.Range("P193").Value = 750
$ws.Range("Q193").Value = 1
$ws.Range("R193").Value = "Hortaliza"

# New row 194.
$ws.Range("A194").Value = $ws.Range("A195").Value()
$ws.Range("B194").Value = $ws.Range("B195").Value()
$ws.Range("C194").Value = $ws.Range("C195").Value()
$ws.Range("D194").Value = 44837
$ws.Range("D194").NumberFormat = $ws.Range("D195").NumberFormat
$ws.Range("E194").Value = $ws.Range("E195").Value()
$ws.Range("F194").Value = $ws.Range("F195").Value()
$ws.Range("G194").Value = $ws.Range("G195").Value()
$ws.Range("H194").Value = $ws.Range("H195").Value()
$ws.Range("I194").Value = "Segunda"
$ws.Range("J194").Value = 150
$ws.Range("K194").Value = 600
$ws.Range("L194").Value = 600
$ws.Range("M194").Value = 600
$ws.Range("N194").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O194").Value = "Provincia de Diguillín"
$ws.Range("P194").Value = 600
$ws.Range("Q194").Value = 1
$ws.Range("R194").Value = "Hortaliza"
